$wb = $excel.ActiveWorkbook

# --- Metrics sheet: update the underlying metric values ---
$metrics = $wb.Worksheets.Item("Metrics")

$metrics.Range("B2").Value = 254992.89
$metrics.Range("B3").Value = 218162.78
$metrics.Range("B4").Value = 78223.7
$metrics.Range("B5").Value = 10459
$metrics.Range("B6").Value = 5457700.0000000009
$metrics.Range("B7").Value = 4618515.74
$metrics.Range("B8").Value = 1610180.5800000003
$metrics.Range("B9").Value = 213166
$metrics.Range("B10").Value = 33923080.989999987
$metrics.Range("B11").Value = 31893790.899999999
$metrics.Range("B12").Value = 11891902.619999995
$metrics.Range("B13").Value = 1310796

# Move the active selection on the Metrics sheet
$metrics.Range("D23").Select() | Out-Null

# --- today sheet: move the active selection (values flow through formulas) ---
$today = $wb.Worksheets.Item("today")
$today.Activate() | Out-Null
$today.Range("G7").Select() | Out-Null
